$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pluralize job titles in A2:A9 (Senior Staff moves to A9, rest stay but become plural),
# and reorder so "Senior Staff" moves from row 2 to row 9 (last before Software Engineer).
$ws.Range("A2").Value = "Archivists"
$ws.Range("A3").Value = "Curators"
$ws.Range("A4").Value = "Data Librarians"
$ws.Range("A5").Value = "Scientists"
$ws.Range("A6").Value = "Policy Specialists"
$ws.Range("A7").Value = "Project Managers"
$ws.Range("A8").Value = "Researchers"
$ws.Range("A9").Value = "Senior Staffs"
$ws.Range("A10").Value = "Software Engineers"

$ws.Range("A6").EntireRow.RowHeight = 26
